$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.789941881126367
$ws.Range("C2").Value = 0.1974320777309231
$ws.Range("D2").Value = 0.1255301292646962
$ws.Range("E2").Value = 0.1219734253762805
$ws.Range("F2").Value = 1.673981462271868
$ws.Range("J2").Value = 0.1521136849595996
$ws.Range("L2").Value = 0.324889485132104
$ws.Range("O2").Value = 4.323744603168905

$ws.Range("B3").Value = 1.65880875767823
$ws.Range("C3").Value = 0.1815793962261409
$ws.Range("D3").Value = 0.1245654263693154
$ws.Range("E3").Value = 0.122817320474268
$ws.Range("F3").Value = 1.686348976962982
$ws.Range("J3").Value = 0.1539056223808393
$ws.Range("L3").Value = 0.3155631487321386
$ws.Range("O3").Value = 4.363243441035053

$ws.Range("B4").Value = 1.5784718402669
$ws.Range("C4").Value = 0.1718055225641706
$ws.Range("D4").Value = 0.1240026283985856
$ws.Range("E4").Value = 0.1233783293661315
$ws.Range("F4").Value = 1.694988679608997
$ws.Range("J4").Value = 0.1550683922271072
$ws.Range("L4").Value = 0.3099255058504014
$ws.Range("O4").Value = 4.390416221686564

$ws.Range("B5").Value = 1.54578101507559
$ws.Range("C5").Value = 0.1678127304105317
$ws.Range("D5").Value = 0.1237807504400266
$ws.Range("E5").Value = 0.1236177360320303
$ws.Range("F5").Value = 1.698772299926361
$ws.Range("J5").Value = 0.1555579589344753
$ws.Range("L5").Value = 0.3076506209180963
$ws.Range("O5").Value = 4.402222931448193

$ws.Range("B6").Value = 1.540355634410218
$ws.Range("C6").Value = 0.1671491427127592
$ws.Range("D6").Value = 0.1237443599418029
$ws.Range("E6").Value = 0.1236581415174598
$ws.Range("F6").Value = 1.699416439484253
$ws.Range("J6").Value = 0.155640201348672
$ws.Range("L6").Value = 0.30727424228823
$ws.Range("O6").Value = 4.404227710792554

$ws.Range("B7").Value = 1.578030766879351
$ws.Range("C7").Value = 0.1717517139373399
$ws.Range("D7").Value = 0.1239996057867998
$ws.Range("E7").Value = 0.1233815143757759
$ws.Range("F7").Value = 1.695038642687372
$ws.Range("J7").Value = 0.1550749309989339
$ws.Range("L7").Value = 0.3098947346134509
$ws.Range("O7").Value = 4.390572481884078

$ws.Range("B8").Value = 1.744691120654011
$ws.Range("C8").Value = 0.191974583074682
$ws.Range("D8").Value = 0.1251913936851778
$ws.Range("E8").Value = 0.1222555151122968
$ws.Range("F8").Value = 1.678028590717297
$ws.Range("J8").Value = 0.1527185760662468
$ws.Range("L8").Value = 0.3216554449716256
$ws.Range("O8").Value = 4.33675722798111

$ws.Range("B9").Value = 2.072863257250674
$ws.Range("C9").Value = 0.2313024441538403
$ws.Range("D9").Value = 0.1277612453177142
$ws.Range("E9").Value = 0.120386797988429
$ws.Range("F9").Value = 1.652979266673412
$ws.Range("J9").Value = 0.148593316301902
$ws.Range("L9").Value = 0.3454160042079479
$ws.Range("O9").Value = 4.254433096844764

$ws.Range("B10").Value = 2.314725130968554
$ws.Range("C10").Value = 0.2599859456025229
$ws.Range("D10").Value = 0.1297893827025689
$ws.Range("E10").Value = 0.119219843425677
$ws.Range("F10").Value = 1.639651993101268
$ws.Range("J10").Value = 0.1458638597767994
$ws.Range("L10").Value = 0.3632917006573138
$ws.Range("O10").Value = 4.208149352101657

$ws.Range("B11").Value = 2.424905070806346
$ws.Range("C11").Value = 0.2729871722924315
$ws.Range("D11").Value = 0.1307420840076574
$ws.Range("E11").Value = 0.1187335069483666
$ws.Range("F11").Value = 1.634694077155117
$ws.Range("J11").Value = 0.144687443680021
$ws.Range("L11").Value = 0.3715135116242294
$ws.Range("O11").Value = 4.19018838638442

$ws.Range("B12").Value = 2.46664813427634
$ws.Range("C12").Value = 0.2779034177767699
$ws.Range("D12").Value = 0.1311071385240439
$ws.Range("E12").Value = 0.1185557308121687
$ws.Range("F12").Value = 1.632975706415749
$ws.Range("J12").Value = 0.1442513359872324
$ws.Range("L12").Value = 0.3746397025579569
$ws.Range("O12").Value = 4.183832819342825

$ws.Range("B13").Value = 2.457657148298381
$ws.Range("C13").Value = 0.2768449336136314
$ws.Range("D13").Value = 0.1310283274593971
$ws.Range("E13").Value = 0.1185937341337624
$ws.Range("F13").Value = 1.633338709386678
$ws.Range("J13").Value = 0.1443448426650766
$ws.Range("L13").Value = 0.3739658563561932
$ws.Range("O13").Value = 4.185181758982196

$ws.Range("B14").Value = 2.42833890255838
$ws.Range("C14").Value = 0.2733917773424537
$ws.Range("D14").Value = 0.1307720315672825
$ws.Range("E14").Value = 0.1187187532266005
$ws.Range("F14").Value = 1.634549516179831
$ws.Range("J14").Value = 0.144651376960713
$ws.Range("L14").Value = 0.371770450349274
$ws.Range("O14").Value = 4.189656568711683

$ws.Range("B15").Value = 2.410383214114688
$ws.Range("C15").Value = 0.2712756949896686
$ws.Range("D15").Value = 0.1306156000754797
$ws.Range("E15").Value = 0.1187961626868592
$ws.Range("F15").Value = 1.63531189448338
$ws.Range("J15").Value = 0.144840358923056
$ws.Range("L15").Value = 0.3704273583797004
$ws.Range("O15").Value = 4.192455612306446

$ws.Range("B16").Value = 2.307527568170542
$ws.Range("C16").Value = 0.2591353148926601
$ws.Range("D16").Value = 0.129727723530479
$ws.Range("E16").Value = 0.1192525213793463
$ws.Range("F16").Value = 1.639998253843217
$ws.Range("J16").Value = 0.1459420535258689
$ws.Range("L16").Value = 0.3627561827599806
$ws.Range("O16").Value = 4.209385481660775

$ws.Range("B17").Value = 2.244467362744558
$ws.Range("C17").Value = 0.2516753413941899
$ws.Range("D17").Value = 0.12919071867627
$ws.Range("E17").Value = 0.1195438755751876
$ws.Range("F17").Value = 1.643156291771746
$ws.Range("J17").Value = 0.1466346108370136
$ws.Range("L17").Value = 0.3580731033043207
$ws.Range("O17").Value = 4.220564490935033

$ws.Range("B18").Value = 2.208211584618539
$ws.Range("C18").Value = 0.2473801501742514
$ws.Range("D18").Value = 0.1288846824027559
$ws.Range("E18").Value = 0.1197156457013637
$ws.Range("F18").Value = 1.645076680384946
$ws.Range("J18").Value = 0.1470390912272075
$ws.Range("L18").Value = 0.3553880081718575
$ws.Range("O18").Value = 4.22728550027503

$ws.Range("B19").Value = 2.195938604262096
$ws.Range("C19").Value = 0.2459251228606263
$ws.Range("D19").Value = 0.1287815518645488
$ws.Range("E19").Value = 0.1197745242891379
$ws.Range("F19").Value = 1.645744740810741
$ws.Range("J19").Value = 0.1471770959702177
$ws.Range("L19").Value = 0.354480344541642
$ws.Range("O19").Value = 4.229611091920759

$ws.Range("B20").Value = 2.251178714349919
$ws.Range("C20").Value = 0.2524699268965946
$ws.Range("D20").Value = 0.1292475906676742
$ws.Range("E20").Value = 0.1195124267725216
$ws.Range("F20").Value = 1.642809351227442
$ws.Range("J20").Value = 0.1465602515186504
$ws.Range("L20").Value = 0.3585707479238494
$ws.Range("O20").Value = 4.219344327380355

$ws.Range("B21").Value = 2.436949847284154
$ws.Range("C21").Value = 0.2744062461765964
$ws.Range("D21").Value = 0.1308471958462505
$ws.Range("E21").Value = 0.1186818587945915
$ws.Range("F21").Value = 1.63418955357146
$ws.Range("J21").Value = 0.1445610859896957
$ws.Range("L21").Value = 0.3724149491910254
$ws.Range("O21").Value = 4.188330100251221

$ws.Range("B22").Value = 2.558479460532965
$ws.Range("C22").Value = 0.288701758751472
$ws.Range("D22").Value = 0.1319176011406142
$ws.Range("E22").Value = 0.1181762684099077
$ws.Range("F22").Value = 1.62948338767022
$ws.Range("J22").Value = 0.1433091618543196
$ws.Range("L22").Value = 0.3815372825010428
$ws.Range("O22").Value = 4.170659783744668

$ws.Range("B23").Value = 2.493606756558336
$ws.Range("C23").Value = 0.2810758279857737
$ws.Range("D23").Value = 0.1313440339592233
$ws.Range("E23").Value = 0.1184427087523918
$ws.Range("F23").Value = 1.631910228084507
$ws.Range("J23").Value = 0.1439723385100038
$ws.Range("L23").Value = 0.3766617784258273
$ws.Range("O23").Value = 4.179852605757759

$ws.Range("B24").Value = 2.248144517197716
$ws.Range("C24").Value = 0.2521107145618089
$ws.Range("D24").Value = 0.1292218704438142
$ws.Range("E24").Value = 0.1195266314849182
$ws.Range("F24").Value = 1.64296587664353
$ws.Range("J24").Value = 0.1465938496794701
$ws.Range("L24").Value = 0.3583457401351779
$ws.Range("O24").Value = 4.21989504743874

$ws.Range("B25").Value = 1.983946642053866
$ws.Range("C25").Value = 0.2206995027477774
$ws.Range("D25").Value = 0.1270413041532308
$ws.Range("E25").Value = 0.1208560968123766
$ws.Range("F25").Value = 1.658865163015754
$ws.Range("J25").Value = 0.1496563324304194
$ws.Range("L25").Value = 0.338914124362887
$ws.Range("O25").Value = 4.274214066712574
